$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = "2021-04-21"
$ws.Range("M3").Value = 250
$ws.Range("N3").Value = 19000
$ws.Range("O3").Value = 20000
$ws.Range("P3").Value = 19500
$ws.Range("Q3").Value = "$/bandeja 18 kilos"
$ws.Range("S3").Value = 1083
$ws.Range("T3").Value = 18

# Row 5
$ws.Range("D5").Value = "2021-05-07"
$ws.Range("M5").Value = 270
$ws.Range("N5").Value = 21000
$ws.Range("O5").Value = 22000
$ws.Range("P5").Value = 21500
$ws.Range("S5").Value = 1194

# Row 7
$ws.Range("D7").Value = "2021-10-20"
$ws.Range("N7").Value = 26000
$ws.Range("O7").Value = 27000
$ws.Range("P7").Value = 26500
$ws.Range("Q7").Value = "$/bandeja 18 kilos"
$ws.Range("S7").Value = 1472
$ws.Range("T7").Value = 18

# Row 8
$ws.Range("D8").Value = "2021-10-18"
$ws.Range("N8").Value = 14000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 14500
$ws.Range("Q8").Value = "$/bandeja 10 kilos"
$ws.Range("S8").Value = 1450
$ws.Range("T8").Value = 10

# Row 9
$ws.Range("D9").Value = "2021-08-10"
$ws.Range("M9").Value = 240
$ws.Range("N9").Value = 10000
$ws.Range("O9").Value = 11000
$ws.Range("P9").Value = 10500
$ws.Range("Q9").Value = "$/bandeja 10 kilos"
$ws.Range("S9").Value = 1050
$ws.Range("T9").Value = 10
